# Insert two new rows above the current row 301 so every existing
# record (rows 301-386) shifts down to rows 303-388, preserving all
# of their data/formatting untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("301:302").Insert()

# Populate the two freshly-inserted rows with the new weekly entries
# (same template columns A,B,C,E,F,G,H,N,O,Q,R as every other row;
# only the date/quality/volume/price columns change).

# Row 301: "Primera" quality entry for 2021-10-10 (serial 44463)
$ws.Range("A301").Value = 8
$ws.Range("B301").Value = "Terminal La Palmera de La Serena"
$ws.Range("C301").Value = "Coquimbo"
$ws.Range("D301").Value = 44463
$ws.Range("E301").Value = 4
$ws.Range("F301").Value = 100112023
$ws.Range("G301").Value = "Brócoli"
$ws.Range("H301").Value = "Sin especificar"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 3460
$ws.Range("K301").Value = 600
$ws.Range("L301").Value = 700
$ws.Range("M301").Value = 650
$ws.Range("N301").Value = "`$/unidad"
$ws.Range("O301").Value = "Provincia del Elquí"
$ws.Range("P301").Value = 650
$ws.Range("Q301").Value = 1
$ws.Range("R301").Value = "Hortaliza"

# Row 302: "Segunda" quality entry for the same date
$ws.Range("A302").Value = 8
$ws.Range("B302").Value = "Terminal La Palmera de La Serena"
$ws.Range("C302").Value = "Coquimbo"
$ws.Range("D302").Value = 44463
$ws.Range("E302").Value = 4
$ws.Range("F302").Value = 100112023
$ws.Range("G302").Value = "Brócoli"
$ws.Range("H302").Value = "Sin especificar"
$ws.Range("I302").Value = "Segunda"
$ws.Range("J302").Value = 1600
$ws.Range("K302").Value = 500
$ws.Range("L302").Value = 550
$ws.Range("M302").Value = 525
$ws.Range("N302").Value = "`$/unidad"
$ws.Range("O302").Value = "Provincia del Elquí"
$ws.Range("P302").Value = 525
$ws.Range("Q302").Value = 1
$ws.Range("R302").Value = "Hortaliza"
